$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.126.41"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.473.41"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "562.14"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "163.58"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.508"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "2.472.80"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("D13").Value = "4.87"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "69.015.66"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.905.82"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "23.75"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "2.461.75"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "10.83"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "339.59"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "7.08"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "3.81"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "1.93"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "67.42"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "3.72"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "2.593.04"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "0.0₃0829"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "439.63"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "1.15"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "156.46"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "19.06"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "0.110"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "17.87"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "0.303"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "4.46"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "37.60"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "1.49"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").Value = "133.69"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "0.0717"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "0.486"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").Value = "0.565"
$ws.Range("E51").Value = "  -0.22%  "
